$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the SVM parameters text (D4): C value changed from 10 to 1
$ws.Range("D4").Value = "C = 1, gamma = 10, kernel = rbf"

# Update the Decision Tree parameters text (D3): max_depth changed from 15 to 14
$ws.Range("D3").Value = "criterion='entropy', max_depth=14, min_samples_leaf=1"

# Narrow column C's width (target stored width 10.6328125 characters;
# the COM width setter snaps to whole-pixel increments, so we pick the
# input value that rounds to the closest representable stored width)
$ws.Columns.Item(3).ColumnWidth = 9.8

# Update the selected cell / active selection
$ws.Range("E6").Select()
